# Update the cryptos Price (D) and Volume(1h) (E) columns with refreshed
# values, as produced by the scheduled GitHub Actions scrape.
#
# Note: several "Price" values look like plain decimal numbers (e.g. "211.99").
# Assigning such a string straight to .Value would make Excel auto-convert it
# to a numeric cell, which would not match the original text-cell layout of
# the sheet. To keep these as text we briefly force a text NumberFormat while
# writing the value, then restore the cell's default ("Normal") style so the
# visible formatting/style of the sheet is left unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "86.994.19"
$ws.Range("E2").Value = "  +3.56%  "
$ws.Range("D3").Value = "3.266.49"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("E4").Value = "  +0.16%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "211.99"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -3.53%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "629.27"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.42%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.376"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +20.91%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.695"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +17.20%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").Value = "3.259.09"
$ws.Range("E10").Value = "  +0.83%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.578"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -5.60%  "
$ws.Range("E12").Value = "  +7.53%  "
$ws.Range("E13").Value = "  -8.33%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "34.24"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +4.48%  "
$ws.Range("D15").Value = "3.867.16"
$ws.Range("E15").Value = "  +1.55%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "5.31"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.27%  "
$ws.Range("D17").Value = "86.887.55"
$ws.Range("E17").Value = "  +3.98%  "
$ws.Range("D18").Value = "3.276.24"
$ws.Range("E18").Value = "  +1.57%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "14.07"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.05%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "3.07"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -6.22%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "433.03"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -3.20%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "8.93"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.50%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.33"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.85%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "7.23"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.34%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "12.51"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +5.00%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "5.13"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("D27").Value = "3.441.20"
$ws.Range("E27").Value = "  +1.38%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "76.12"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -3.15%  "
$ws.Range("E29").Value = "  +3.78%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("E31").Value = "  +15.24%  "
$ws.Range("E32").Value = "  +0.47%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "8.82"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -4.63%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "545.67"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -4.79%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.43"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -4.56%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.94"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.88%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "6.93"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +10.96%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.138"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -10.43%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "22.52"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -3.33%  "
$ws.Range("E40").Value = "  -0.06%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "21.57"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +3.04%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.396"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.33%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.00"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -3.56%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.93"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.60%  "
$ws.Range("E45").Value = "  +0.02%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "157.33"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.72%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "179.45"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -5.22%  "
$ws.Range("E48").Value = "  -1.29%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.31"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.81%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "4.25"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.14%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.627"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.88%  "
